$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.939.60"
$ws.Range("E2").Value = "  +6.15%  "
$ws.Range("D3").Value = "3.542.67"
$ws.Range("E3").Value = "  +9.82%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'565.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.26%  "
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "'189.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.24%  "
$ws.Range("D7").Value = "3.537.35"
$ws.Range("E7").Value = "  +9.65%  "
$ws.Range("D8").Value = "'0.618"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.98%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'0.636"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.47%  "
$ws.Range("E11").Value = "  +14.19%  "
$ws.Range("D12").Value = "'55.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.03%  "
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.68%  "
$ws.Range("D14").Value = "'9.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.64%  "
$ws.Range("D15").Value = "4.108.31"
$ws.Range("E15").Value = "  +9.79%  "
$ws.Range("D16").Value = "3.540.93"
$ws.Range("E16").Value = "  +9.64%  "
$ws.Range("D17").Value = "'0.122"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.91%  "
$ws.Range("D18").Value = "66.985.66"
$ws.Range("E18").Value = "  +6.39%  "
$ws.Range("D19").Value = "'18.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.82%  "
$ws.Range("D20").Value = "'12.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.05%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.56%  "
$ws.Range("D22").Value = "'430.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +17.76%  "
$ws.Range("D23").Value = "'4.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.07%  "
$ws.Range("D24").Value = "'85.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.10%  "
$ws.Range("D25").Value = "'4.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.10%  "
$ws.Range("D26").Value = "'11.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").Value = "'2.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.92%  "
$ws.Range("D28").Value = "'6.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "'12.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.67%  "
$ws.Range("D30").Value = "'9.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.29%  "
$ws.Range("D31").Value = "'30.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.21%  "
$ws.Range("D32").Value = "'645.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").Value = "'6.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.52%  "
$ws.Range("D34").Value = "'11.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.77%  "
$ws.Range("D35").Value = "'0.113"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.94%  "
$ws.Range("D36").Value = "'59.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.28%  "
$ws.Range("D37").Value = "'0.151"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +23.67%  "
$ws.Range("D38").Value = "0.0₃0825"
$ws.Range("E38").Value = "  +16.41%  "
$ws.Range("D39").Value = "'38.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.06%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'0.393"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.55%  "
$ws.Range("D42").Value = "'3.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.09%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "3.058.31"
$ws.Range("E44").Value = "  +6.75%  "
$ws.Range("E45").Value = "  +5.17%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'2.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.71%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.51%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0421"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.42%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'2.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.80%  "
$ws.Range("E50").Value = "  +5.80%  "
$ws.Range("D51").Value = "'8.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.83%  "
